$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: rename the "U" / extra "F" header labels ---
# E15 (H2F) and F15 (H4F) stay the same.
# G15: H3U -> H6F
# H15: H4U -> H8U
# I15: H5U -> H9U
# J15 (new): H10U
$ws.Range("G15").Value = "H6F"
$ws.Range("H15").Value = "H8U"
$ws.Range("I15").Value = "H9U"
$ws.Range("J15").Value = "H10U"

# --- Row 14: add a new header cell "ctr" in column J ---
$ws.Range("J14").Value = "ctr"

# --- Row 16: add a formula in column J referencing C16 ---
$ws.Range("J16").Formula = "=C16"

# --- Update the selected range shown in the sheet view ---
$ws.Range("A16").Select() | Out-Null
